$d = $word.ActiveDocument

# Locate the target paragraph: "This is a syntax guide to the Arduino Language...doubles. "
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("This is a syntax guide")) {
        $target = $p
        break
    }
}

$start = $target.Range.Start

# A hidden "_GoBack" bookmark sits immediately after "doubles" (and right before the
# trailing ". " that closes the paragraph) - locate it by searching for that word so we
# don't depend on hard-coded offsets.
$searchRange = $target.Range.Duplicate
$searchRange.Find.Execute("doubles", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bookmarkPos = $searchRange.End

$newText = "This is a syntax guide to the Arduino Language (AL). There are several aspects to the Arduino Language that are different than the conventional programming language: the core being the language is centered on reading like English and is loosely-typed. AL is an imperative and procedural programming language. The intention of this language is to educate people who desire to learn programming. The ease and flexibility of the language allows for a dynamic experience along with a sense of learning."

# Replace everything from the paragraph start through "doubles" (i.e. right up to, but not
# including, the bookmark) with the full new sentence in one shot - this makes the bookmark
# shift to follow the freshly inserted text instead of getting deleted or stranded.
$r = $d.Range($start, $bookmarkPos)
$r.Text = $newText

# Remove the old trailing ". " run that used to follow the bookmark - nothing should follow
# the bookmark now except the paragraph mark.
$endOfText = $start + $newText.Length
$paraEnd = $target.Range.End
$needsTrim = $paraEnd - 1 -gt $endOfText
if ($needsTrim) {
    $trailing = $d.Range($endOfText, $paraEnd - 1)
    $trailing.Text = ""
}
